# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Update the curated-dimension metadata rows of the report sheet:
#  - Row 2 (sdmx/iaest dimension identifiers):
#      C2: iaest-measure:municipio-nombre            -> sdmx-dimension:refArea
#      E2: iaest-dimension:porcentaje-actividades-desarrollo-rural
#                                                      -> iaest-measure:porcentaje-actividades-desarrollo-rural
#  - Row 3 (dim/medida classification):
#      C3: medida -> dim
#      E3: dim    -> medida
#  - Row 4 (RDF type / URI template):
#      C4: xsd:int      -> URI-Municipio
#      E4: skos:Concept -> xsd:int
#  - Row 5 (mapping workbook references):
#      E5: mapping-porcentaje-actividades-desarrollo-rural.xlsx -> (removed entirely)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "iaest-measure:porcentaje-actividades-desarrollo-rural"

$ws.Range("C3").Value = "dim"
$ws.Range("E3").Value = "medida"

$ws.Range("C4").Value = "URI-Municipio"
$ws.Range("E4").Value = "xsd:int"

# Fully clear E5 (content + formatting) so the cell element is dropped, not
# just emptied, matching the removal of the mapping-file reference.
$ws.Range("E5").Clear()
